$d = $word.ActiveDocument

$replacements = @(
    @("59×25=", "98×99="),
    @("84×51=", "69×63="),
    @("66×84=", "38×70="),
    @("90×77=", "95×52="),
    @("33×50=", "28×64="),
    @("71×25=", "27×47="),
    @("72×77=", "24×73="),
    @("44×87=", "29×41="),
    @("87×34=", "48×12="),
    @("36×98=", "81×62="),
    @("53×88=", "28×68="),
    @("77×67=", "45×14="),
    @("61×54=", "14×97="),
    @("30×54=", "22×80="),
    @("72×18=", "45×17="),
    @("83×44=", "41×34="),
    @("97×65=", "96×84="),
    @("95×30=", "98×34="),
    @("86×62=", "57×95="),
    @("14×14=", "89×30="),
    @("99×85=", "75×85="),
    @("69×68=", "76×86="),
    @("16×16=", "97×22="),
    @("91×93=", "71×88="),
    @("81×63=", "28×29=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
